# Generate Report for Handoff
# Replace the old localization-status GUID / content-hash used throughout the
# workbook with the newly generated ones, and bump the two "Latest Handoff
# Datetime" timestamps to reflect the new handoff run.

$oldGuid = "3f9d76b3-ebc3-4373-b9ae-81d2bfcf476f"
$newGuid = "8a6eb94e-dc8f-4a69-88a7-a0f9a48906a9"
$oldHash = "555cb7280b876f3270c80da6bc0f4b1d28b65375"
$newHash = "86847819c0b0992d8e5759dfe6448c7d18e83ef1"

$wb = $excel.ActiveWorkbook

$sheetOverview = $wb.Worksheets.Item("Overview")
$sheetZhCn     = $wb.Worksheets.Item("zh-cn")
$sheetDeDe     = $wb.Worksheets.Item("de-de")

# ---- new display text / cell text -----------------------------------------
$newMdName    = "$newGuid.md"
$newZhXlfName = "$newGuid.$newHash.zh-cn.xlf"
$newDeXlfName = "$newGuid.$newHash.de-de.xlf"
$newZhDateTime = "2016-02-29 04:19:43"
$newDeDateTime = "2016-02-29 04:19:58"

# ---- hyperlink target addresses (unchanged from the original workbook) ----
$mdAddress    = "https://github.com/OpenLocalizationTest/oltest/blob/b9c6accbedc5f45413f0715f22b900c13a2b5c36/e2e/$oldGuid.md"
$configAddress = "https://github.com/OpenLocalizationTest/oltest/blob/b9c6accbedc5f45413f0715f22b900c13a2b5c36/.localization-config"
$zhXlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7c9a9eb0d5b8882c1b062375ecc664bb6330633e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$oldGuid.$oldHash.zh-cn.xlf"
$deXlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/86ad12d516cacd45a2191113eae78e9d840f8e5b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$oldGuid.$oldHash.de-de.xlf"

function Restore-HyperlinkLook($range) {
    # Hyperlinks.Add() re-styles the cell with Excel's built-in themed
    # "Hyperlink" cell style; put back the workbook's original custom look
    # (underlined, RGB 0x6495ED) so the visual style is unchanged.
    $range.Font.Underline = $true
    $range.Font.Color = 15570276   # BGR-packed 0x6495ED (matches original theme)
}

# ---- Overview sheet ---------------------------------------------------------
$sheetOverview.Hyperlinks.Delete()
$sheetOverview.Hyperlinks.Add($sheetOverview.Range("A2"), $mdAddress, "", "", $newMdName)
$sheetOverview.Hyperlinks.Add($sheetOverview.Range("A3"), $configAddress, "", "", ".localization-config")
Restore-HyperlinkLook($sheetOverview.Range("A2"))
Restore-HyperlinkLook($sheetOverview.Range("A3"))

# ---- zh-cn sheet ------------------------------------------------------------
$sheetZhCn.Hyperlinks.Delete()
$sheetZhCn.Hyperlinks.Add($sheetZhCn.Range("A2"), $mdAddress, "", "", $newMdName)
$sheetZhCn.Hyperlinks.Add($sheetZhCn.Range("C2"), $zhXlfAddress, "", "", $newZhXlfName)
$sheetZhCn.Hyperlinks.Add($sheetZhCn.Range("A3"), $configAddress, "", "", ".localization-config")
$sheetZhCn.Range("D2").Value = $newZhDateTime
Restore-HyperlinkLook($sheetZhCn.Range("A2"))
Restore-HyperlinkLook($sheetZhCn.Range("C2"))
Restore-HyperlinkLook($sheetZhCn.Range("A3"))
# Writing a literal string into D2 can make the engine drop the shared
# "yyyy-mm-dd HH:mm:ss" number format from its style slot (and, since D3
# shares that same slot, D3 would silently lose it too) - pin it back on
# both cells explicitly.
$sheetZhCn.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$sheetZhCn.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---- de-de sheet ------------------------------------------------------------
$sheetDeDe.Hyperlinks.Delete()
$sheetDeDe.Hyperlinks.Add($sheetDeDe.Range("A2"), $mdAddress, "", "", $newMdName)
$sheetDeDe.Hyperlinks.Add($sheetDeDe.Range("C2"), $deXlfAddress, "", "", $newDeXlfName)
$sheetDeDe.Hyperlinks.Add($sheetDeDe.Range("A3"), $configAddress, "", "", ".localization-config")
$sheetDeDe.Range("D2").Value = $newDeDateTime
Restore-HyperlinkLook($sheetDeDe.Range("A2"))
Restore-HyperlinkLook($sheetDeDe.Range("C2"))
Restore-HyperlinkLook($sheetDeDe.Range("A3"))
$sheetDeDe.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$sheetDeDe.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

Write-Output "done"
